$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.104.41"
$ws.Range("E2").Value = "  -2.33%  "
$ws.Range("D3").Value = "1.849.71"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "0.6936"
$ws.Range("E5").Value = "  -5.22%  "
$ws.Range("D6").Value = "237.75"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("D7").Value = "0.9986"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "0.07725"
$ws.Range("E8").Value = "  +8.32%  "
$ws.Range("D9").Value = "0.3043"
$ws.Range("E9").Value = "  -3.06%  "
$ws.Range("D10").Value = "23.29"
$ws.Range("E10").Value = "  -4.71%  "
$ws.Range("D11").Value = "0.08115"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "1.851.54"
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("D13").Value = "0.7253"
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("D14").Value = "5.209"
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").Value = "89.03"
$ws.Range("E15").Value = "  -3.70%  "
$ws.Range("D16").Value = "29.106.94"
$ws.Range("E16").Value = "  -2.28%  "
$ws.Range("D17").Value = "0.000007837"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "5.740"
$ws.Range("E18").Value = "  -4.42%  "
$ws.Range("D19").Value = "13.19"
$ws.Range("E19").Value = "  -1.38%  "
$ws.Range("D20").Value = "236.25"
$ws.Range("E20").Value = "  -4.71%  "
$ws.Range("D21").Value = "0.9990"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "2.100.30"
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("D23").Value = "0.9993"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "7.597"
$ws.Range("E24").Value = "  -1.78%  "
$ws.Range("D25").Value = "8.982"
$ws.Range("D26").Value = "161.08"
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("D27").Value = "0.1429"
$ws.Range("E27").Value = "  -7.14%  "
$ws.Range("D28").Value = "18.06"
$ws.Range("E28").Value = "  -2.61%  "
$ws.Range("D29").Value = "1.973"
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("D30").Value = "1.395"
$ws.Range("E30").Value = "  -3.69%  "
$ws.Range("D31").Value = "4.500"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").Value = "1.487"
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("D33").Value = "4.008"
$ws.Range("E33").Value = "  -4.17%  "
$ws.Range("D34").Value = "0.05227"
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("D35").Value = "1.180"
$ws.Range("E35").Value = "  -4.09%  "
$ws.Range("D36").Value = "0.7053"
$ws.Range("E36").Value = "  -4.69%  "
$ws.Range("D37").Value = "1.021"
$ws.Range("E37").Value = "  +2.32%  "
$ws.Range("D38").Value = "2.644"
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("D39").Value = "0.01854"
$ws.Range("E39").Value = "  -4.26%  "
$ws.Range("D40").Value = "2.670"
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("D41").Value = "0.9109"
$ws.Range("E41").Value = "  +4.79%  "
$ws.Range("D42").Value = "1.096.44"
$ws.Range("E42").Value = "  +4.79%  "
$ws.Range("D43").Value = "6.006"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").Value = "0.4269"
$ws.Range("E44").Value = "  -4.50%  "
$ws.Range("D45").Value = "70.62"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").Value = "0.9982"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").Value = "102.86"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").Value = "1.765"
$ws.Range("E48").Value = "  -2.82%  "
$ws.Range("D49").Value = "1.996.12"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").Value = "9.157"
$ws.Range("E50").Value = "  -4.08%  "
$ws.Range("D51").Value = "6.974"
$ws.Range("E51").Value = "  -6.14%  "
